# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for the zh-cn and de-de handback rows, marks the items as handed back, and widens
# the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a3a33c9560ecb3280b64efe918e32082f81decd/e2e/a.md"
$hyperlinkColor = 15570276   # OLE (BGR) encoding of RGB 64,95,ED -> matches the workbook's existing HyperLink font color

function Set-HandbackRow {
    param($ws, $row, $targetFileName)

    # Column I = Latest Target File -> "a.md", with a hyperlink just like column A's a.md link
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value2 = "a.md"
    $ws.Hyperlinks.Add($iCell, $aMdUrl, [type]::Missing, [type]::Missing, "a.md")
    $iCell.Font.Underline = 2
    $iCell.Font.Color = $hyperlinkColor

    # Column J = Latest Handback File
    $ws.Cells.Item($row, 10).Value2 = $targetFileName

    # Column K = Latest Handback DateTime
    $ws.Cells.Item($row, 11).Value2 = "2016-08-29 08:39:41"
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2:C3").Value2 = "Handed back: in sync with en-US"
Set-HandbackRow $wsZh 2 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Set-HandbackRow $wsZh 3 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Cells.Item(2, 11).Value2 = "2016-08-29 08:39:41"
$wsZh.Cells.Item(3, 11).Value2 = "2016-08-29 08:39:41"
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2:C3").Value2 = "Handed back: in sync with en-US"
Set-HandbackRow $wsDe 2 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Set-HandbackRow $wsDe 3 "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Cells.Item(2, 11).Value2 = "2016-08-29 08:39:48"
$wsDe.Cells.Item(3, 11).Value2 = "2016-08-29 08:39:48"
$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1

# ---- Overview sheet: widen the language-status columns to match ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

Write-Output "Handback report generated"
